# Updates the cryptocurrency price ("D") and 1h volume-change ("E") columns
# to reflect the refreshed values captured by the scheduled data pull.
# Price cells that look like plain numbers (e.g. "227.74") are written via
# Formula with a leading apostrophe so Excel keeps them as text (matching the
# original inline-string cell type) instead of silently converting them to
# numeric values. Prices containing extra "." separators (e.g. "37.810.47")
# are already unambiguous text, so a plain Value assignment is enough.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.810.47'
$ws.Range("E2").Value = '  -0.06%  '
$ws.Range("D3").Value = '2.045.88'
$ws.Range("E3").Value = '  +0.70%  '
$ws.Range("D5").Formula = '''227.74'
$ws.Range("E5").Value = '  +0.10%  '
$ws.Range("E6").Value = '  -0.54%  '
$ws.Range("D7").Formula = '''59.70'
$ws.Range("E7").Value = '  +0.21%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("E9").Value = '  -2.04%  '
$ws.Range("D10").Formula = '''0.0834'
$ws.Range("E10").Value = '  +2.68%  '
$ws.Range("E11").Value = '  -0.05%  '
$ws.Range("D12").Value = '2.349.89'
$ws.Range("E12").Value = '  +0.77%  '
$ws.Range("E13").Value = '  -1.53%  '
$ws.Range("D14").Formula = '''21.41'
$ws.Range("E14").Value = '  +1.11%  '
$ws.Range("E15").Value = '  +6.08%  '
$ws.Range("D16").Formula = '''0.763'
$ws.Range("E16").Value = '  +0.14%  '
$ws.Range("D17").Value = '2.035.72'
$ws.Range("E17").Value = '  -0.47%  '
$ws.Range("D18").Value = '37.773.19'
$ws.Range("E18").Value = '  -0.02%  '
$ws.Range("D19").Formula = '''69.49'
$ws.Range("E19").Value = '  -0.68%  '
$ws.Range("E20").Value = '  -1.84%  '
$ws.Range("D21").Value = '0.0₃0828'
$ws.Range("E21").Value = '  +0.47%  '
$ws.Range("D22").Formula = '''222.46'
$ws.Range("E22").Value = '  -1.09%  '
$ws.Range("E24").Value = '  +0.78%  '
$ws.Range("D26").Formula = '''169.01'
$ws.Range("E26").Value = '  +2.48%  '
$ws.Range("E27").Value = '  +0.53%  '
$ws.Range("E28").Value = '  -1.04%  '
$ws.Range("E29").Value = '  -0.93%  '
$ws.Range("D30").Formula = '''1.29'
$ws.Range("E30").Value = '  +0.10%  '
$ws.Range("E31").Value = '  -0.54%  '
$ws.Range("E32").Value = '  +8.03%  '
$ws.Range("D33").Formula = '''4.37'
$ws.Range("E33").Value = '  -1.40%  '
$ws.Range("D34").Formula = '''4.53'
$ws.Range("E34").Value = '  +0.57%  '
$ws.Range("E35").Value = '  +0.14%  '
$ws.Range("D36").Formula = '''6.47'
$ws.Range("E36").Value = '  +1.83%  '
$ws.Range("E37").Value = '  +3.97%  '
$ws.Range("D38").Formula = '''3.48'
$ws.Range("E38").Value = '  +7.14%  '
$ws.Range("D39").Formula = '''1.00'
$ws.Range("E39").Value = '  -0.06%  '
$ws.Range("D40").Formula = '''18.40'
$ws.Range("E40").Value = '  +9.35%  '
$ws.Range("D41").Value = '1.524.30'
$ws.Range("E41").Value = '  +0.04%  '
$ws.Range("D42").Formula = '''97.53'
$ws.Range("E42").Value = '  +0.81%  '
$ws.Range("D43").Formula = '''0.0215'
$ws.Range("E43").Value = '  -1.74%  '
$ws.Range("E44").Value = '  -0.06%  '
$ws.Range("D45").Formula = '''4.20'
$ws.Range("E45").Value = '  +3.48%  '
$ws.Range("E46").Value = '  -3.06%  '
$ws.Range("E47").Value = '  -0.36%  '
$ws.Range("E48").Value = '  -0.20%  '
$ws.Range("E50").Value = '  +0.16%  '
$ws.Range("D51").Value = '2.238.42'
$ws.Range("E51").Value = '  +0.82%  '
